# Fruta / hortaliza, semanal
#
# A new weekly price record is added for "Terminal La Palmera de La Serena -
# Alcachofa". The previous latest row (18) is pushed down to become row 19
# (unchanged, historical data), and row 18 is populated with this week's
# updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting the existing row 18 (and its
# formatting) down to row 19 — this preserves the historical record exactly.
$ws.Rows.Item(18).Insert(-4121)   # xlShiftDown

# Populate the now-empty row 18 with this week's new record.
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44858
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 100112013
$ws.Range("G18").Value = "Alcachofa"
$ws.Range("H18").Value = "Española"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 9500
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9750
$ws.Range("N18").Value = "$/caja 30 unidades"
$ws.Range("O18").Value = "Provincia del Elquí"
$ws.Range("P18").Value = 325
$ws.Range("Q18").Value = 30
$ws.Range("R18").Value = "Hortaliza"
